$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextCell 'D2' '317.70'
Set-TextCell 'E2' '3.76%'
Set-TextCell 'D3' '39.84'
Set-TextCell 'E3' '2.55%'
Set-TextCell 'D4' '5.153'
Set-TextCell 'E4' '0.75%'
Set-TextCell 'D5' '0.08230'
Set-TextCell 'E5' '1.86%'
Set-TextCell 'D6' '2.059'
Set-TextCell 'E6' '6.31%'
Set-TextCell 'D7' '8.352'
Set-TextCell 'E7' '4.30%'
Set-TextCell 'B8' 'GateToken'
Set-TextCell 'C8' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell 'D8' '4.328'
Set-TextCell 'E8' '3.34%'
Set-TextCell 'B9' 'MXToken'
Set-TextCell 'C9' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D9' '0.9390'
Set-TextCell 'E9' '0.87%'
Set-TextCell 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextCell 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 'D10' '0.1373'
Set-TextCell 'E10' '-6.17%'
Set-TextCell 'B11' 'WazirX'
Set-TextCell 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 'D11' '0.1984'
Set-TextCell 'E11' '3.04%'
Set-TextCell 'B12' 'MandalaExchangeToken'
Set-TextCell 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 'D12' '0.09116'
Set-TextCell 'E12' '0.80%'
Set-TextCell 'B13' 'BitrueCoin'
Set-TextCell 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 'D13' '0.03512'
Set-TextCell 'E13' '0.07%'
Set-TextCell 'B14' 'BitMartToken'
Set-TextCell 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 'D14' '0.09803'
Set-TextCell 'E14' '-0.01%'
Set-TextCell 'B15' 'BitForexToken'
Set-TextCell 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 'D15' '0.001405'
Set-TextCell 'E15' '0.39%'
Set-TextCell 'B16' 'TigerCash'
Set-TextCell 'C16' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell 'D16' '0.006136'
Set-TextCell 'E16' '2.64%'
Set-TextCell 'B17' 'LEO'
Set-TextCell 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 'D17' '3.692'
Set-TextCell 'E17' '-2.32%'
Set-TextCell 'D18' '3.328'
Set-TextCell 'E18' '-3.56%'
Set-TextCell 'D19' '0.3497'
Set-TextCell 'E19' '1.07%'
Set-TextCell 'E20' '-2.83%'
Set-TextCell 'D21' '4.962'
Set-TextCell 'E21' '5.99%'
Set-TextCell 'D22' '0.2448'
Set-TextCell 'E22' '1.37%'
Set-TextCell 'D23' '0.04355'
Set-TextCell 'E23' '-0.31%'
Set-TextCell 'E24' '-0.80%'
Set-TextCell 'D25' '0.004829'
Set-TextCell 'E25' '12.82%'
Set-TextCell 'E26' '-0.25%'
Set-TextCell 'D27' '0.0003995'
Set-TextCell 'E27' '-10.17%'
Set-TextCell 'E39' '9.05%'
Set-TextCell 'D40' '0.05223'
Set-TextCell 'E40' '3.08%'
Set-TextCell 'D41' '0.007756'
Set-TextCell 'E41' '1.96%'
Set-TextCell 'D42' '0.009673'
Set-TextCell 'E42' '-0.64%'
Set-TextCell 'D43' '0.1406'
Set-TextCell 'E43' '4.16%'
Set-TextCell 'D44' '0.002048'
Set-TextCell 'E44' '-3.62%'
Set-TextCell 'D45' '0.009659'
Set-TextCell 'E45' '-2.67%'
Set-TextCell 'D46' '0.00006607'
Set-TextCell 'E46' '6.27%'
Set-TextCell 'D47' '0.00000000749'
Set-TextCell 'E47' '-0.35%'
Set-TextCell 'B48' 'BOLO'
Set-TextCell 'C48' 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextCell 'D48' '0.002940'
Set-TextCell 'E48' '1.70%'
Set-TextCell 'B49' 'CoinbaseStockToken'
Set-TextCell 'C49' 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextCell 'D49' '0.001688'
Set-TextCell 'E49' '-6.44%'
Set-TextCell 'D50' '0.00002098'
Set-TextCell 'E50' '-0.35%'
Set-TextCell 'D51' '0.0001998'
Set-TextCell 'E51' '-0.35%'
